$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D/B/C columns keep their original text representation (avoid Excel
# auto-coercing numeric-looking strings like "1.00" or "0.0920" into numbers,
# which would silently drop significant trailing zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.782.87'
$ws.Range("E2").Value = '  -1.15%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.029.85'
$ws.Range("E3").Value = '  -2.03%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.36'
$ws.Range("E5").Value = '  -1.76%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.613'
$ws.Range("E6").Value = '  -0.72%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.95'
$ws.Range("E7").Value = '  +3.11%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.386'
$ws.Range("E9").Value = '  -0.45%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0817'
$ws.Range("E10").Value = '  +0.96%  '

$ws.Range("E11").Value = '  +0.06%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.61'
$ws.Range("E12").Value = '  -0.19%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.328.36'
$ws.Range("E13").Value = '  -2.08%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.12'
$ws.Range("E14").Value = '  +1.49%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.754'
$ws.Range("E15").Value = '  +0.05%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.22'
$ws.Range("E16").Value = '  -0.80%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.039.12'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.774.60'
$ws.Range("E18").Value = '  -0.87%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.01'
$ws.Range("E19").Value = '  -3.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.75'
$ws.Range("E20").Value = '  -0.44%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0824'
$ws.Range("E21").Value = '  -1.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '225.59'
$ws.Range("E22").Value = '  +0.09%  '

$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("E24").Value = '  -2.03%  '

$ws.Range("E25").Value = '  -2.35%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.24'
$ws.Range("E26").Value = '  -1.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.79'
$ws.Range("E27").Value = '  -0.72%  '

$ws.Range("E28").Value = '  -3.56%  '

$ws.Range("E29").Value = '  -1.44%  '

$ws.Range("E30").Value = '  -6.43%  '

$ws.Range("E31").Value = '  +1.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.43'
$ws.Range("E32").Value = '  -3.11%  '

$ws.Range("E33").Value = '  +3.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0602'
$ws.Range("E34").Value = '  -2.43%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.47'
$ws.Range("E35").Value = '  -3.09%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.37'
$ws.Range("E36").Value = '  +5.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.25'
$ws.Range("E37").Value = '  -5.88%  '

$ws.Range("E38").Value = '  -1.79%  '

$ws.Range("E39").Value = '  -0.17%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.536.32'
$ws.Range("E40").Value = '  +3.51%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0217'
$ws.Range("E41").Value = '  -1.10%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '96.74'
$ws.Range("E42").Value = '  -1.90%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.63'
$ws.Range("E43").Value = '  -1.29%  '

$ws.Range("E44").Value = '  -1.88%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0920'
$ws.Range("E45").Value = '  -3.36%  '

$ws.Range("E46").Value = '  -2.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.97'
$ws.Range("E47").Value = '  -1.97%  '

$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  -2.51%  '

$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.96'
$ws.Range("E49").Value = '  -0.42%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.13'
$ws.Range("E50").Value = '  -0.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.218.70'
$ws.Range("E51").Value = '  -1.88%  '
